# The document contains one table with a "stuklijst" (parts list). A number
# of cells contain English/technical terms that Word's spell checker flags;
# this edit wraps those terms in <w:proofErr w:type="spellStart"/> ...
# <w:proofErr w:type="spellEnd"/> markers (as Word does automatically while
# a user types/reviews a document), splitting the affected run(s) so the
# flagged word is its own run. A trailing empty paragraph is also added
# after the table.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rNs = 'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

function Set-CellSpellSplit($row, $col, $beforeText, $word1, $afterText) {
    $cell = $tbl.Cell($row, $col)
    $range = $cell.Range

    $xml = "<w:p $wNs>"
    if ($beforeText -ne "") {
        $xml += "<w:r><w:t xml:space=`"preserve`">$beforeText</w:t></w:r>"
    }
    $xml += '<w:proofErr w:type="spellStart"/>'
    $xml += "<w:r><w:t>$word1</w:t></w:r>"
    $xml += '<w:proofErr w:type="spellEnd"/>'
    if ($afterText -ne "") {
        $xml += "<w:r><w:t xml:space=`"preserve`">$afterText</w:t></w:r>"
    }
    $xml += "</w:p>"

    $range.InsertXML($xml)
}

function Set-HyperlinkCellSpellSplit($row, $col, $relId, $beforeText, $word1, $afterText) {
    $cell = $tbl.Cell($row, $col)
    $range = $cell.Range

    $xml = "<w:p $wNs $rNs>"
    $xml += "<w:hyperlink r:id=`"$relId`" w:tgtFrame=`"_blank`" w:history=`"1`">"
    if ($beforeText -ne "") {
        $xml += "<w:r><w:t xml:space=`"preserve`">$beforeText</w:t></w:r>"
    }
    $xml += '<w:proofErr w:type="spellStart"/>'
    $xml += "<w:r><w:t>$word1</w:t></w:r>"
    $xml += '<w:proofErr w:type="spellEnd"/>'
    if ($afterText -ne "") {
        $xml += "<w:r><w:t xml:space=`"preserve`">$afterText</w:t></w:r>"
    }
    $xml += "</w:hyperlink>"
    $xml += "</w:p>"

    $range.InsertXML($xml)
}

# Locate the rows by their current label text so the edit is resilient to
# row ordering.
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $label = $tbl.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13)

    if ($label -eq "NI MyRio") {
        Set-CellSpellSplit $i 1 "NI " "MyRio" ""
    }
    elseif ($label -eq "Tiny breadboard") {
        Set-CellSpellSplit $i 1 "Tiny " "breadboard" ""
    }
    elseif ($label -eq "Ball Caster") {
        Set-HyperlinkCellSpellSplit $i 1 "rId7" "Ball " "Caster" ""
    }
    elseif ($label -eq "Micro metal gear motor beugel") {
        Set-HyperlinkCellSpellSplit $i 1 "rId9" "Micro metal " "gear" " motor beugel"
    }
}

# Insert an extra empty paragraph right after the table (before the existing
# trailing empty paragraph / sectPr). NOTE: $d.Paragraphs becomes stale after
# the InsertXML calls above (this runtime doesn't refresh paragraph ranges
# once a raw-XML mutation has happened), so recompute the insertion point
# from $d.Content.End instead of using $d.Paragraphs.Last.
$endPos = $d.Content.End
$insertPoint = $d.Range($endPos - 1, $endPos - 1)
$insertPoint.InsertXML("<w:p $wNs/>")
